$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '290.51'
    'E2' = '-3.80%'
    'D3' = '30.74'
    'E3' = '-5.97%'
    'D4' = '4.954'
    'D5' = '0.07212'
    'E5' = '-7.82%'
    'D6' = '1.786'
    'E6' = '-8.83%'
    'D7' = '7.671'
    'E7' = '-2.41%'
    'D8' = '3.723'
    'D9' = '0.8962'
    'E9' = '-3.34%'
    'D10' = '0.1654'
    'E10' = '-6.52%'
    'D11' = '0.07708'
    'E11' = '-1.73%'
    'D12' = '0.08032'
    'E12' = '-7.10%'
    'D13' = '0.03049'
    'E13' = '-3.93%'
    'D14' = '0.1002'
    'E14' = '0.01%'
    'D15' = '0.001494'
    'E15' = '-1.32%'
    'D16' = '0.005660'
    'E16' = '-1.33%'
    'D17' = '3.476'
    'E17' = '0.36%'
    'D18' = '2.084'
    'E18' = '-3.26%'
    'D19' = '0.3315'
    'E19' = '-0.44%'
    'D20' = '0.1305'
    'E20' = '-1.04%'
    'E21' = '-6.57%'
    'E22' = '5.42%'
    'D23' = '0.04522'
    'E23' = '-1.14%'
    'E24' = '-0.92%'
    'D25' = '0.004018'
    'E25' = '-9.43%'
    'D26' = '0.0001250'
    'E26' = '-0.15%'
    'D39' = '0.01598'
    'E39' = '-5.98%'
    'D40' = '0.04392'
    'E40' = '-7.50%'
    'D41' = '0.007333'
    'E41' = '-5.34%'
    'E42' = '-3.43%'
    'D43' = '0.007692'
    'D44' = '0.002060'
    'E44' = '-12.10%'
    'D45' = '0.009215'
    'E45' = '-21.25%'
    'D46' = '0.00005845'
    'E46' = '-6.52%'
    'E47' = '-0.15%'
    'D48' = '2.246'
    'E48' = '173.66%'
    'D49' = '0.003000'
    'E49' = '-3.34%'
    'E50' = '-0.15%'
    'E51' = '-0.15%'
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
